$d = $word.ActiveDocument

# --- Step 1: remove the existing "_GoBack" bookmark from its current spot
# (right after "...triggered or cycling operation. " and before the trailing
# two-space run). It gets re-added later, at the end of the new paragraph.
$bm = $d.Bookmarks.Item("_GoBack")
$bm.Delete()

# --- Step 2: find the paragraph that currently just says "Text here" and
# replace it with the new descriptive copy about the Dual Random Generator.
$dash = [char]0x2013
$newText = "The Dual Random Generator device can be used to obtain many flavors of random modulation. Each side of the Dual Random Generator can be set to one of three modes " + $dash + " The first mode outputs random stepped voltages, the second mode outputs random smoothed voltages and the last mode outputs random gates or timing pulses. Each mode has controls for rate, variance and scaling functions like offset and depth. The Dual Random Generator "

$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $ptext = $p.Range.Text.TrimEnd("`r")
    if ($ptext -eq "Text here") {
        $r = $p.Range
        # drop the trailing paragraph mark from the replace range
        $null = $r.MoveEnd(1, -1)
        # Append a one-character placeholder ("Z") so the insertion point for
        # the bookmark below is not exactly the paragraph-end position.
        $r.Text = $newText + "Z"

        $p2 = $d.Paragraphs.Item($i)
        $paraEnd = $p2.Range.End
        $bmPos = $paraEnd - 2  # right before the "Z" placeholder

        $bmRange = $d.Range($bmPos, $bmPos)
        $d.Bookmarks.Add("_GoBack", $bmRange)

        # remove the "Z" placeholder now that the bookmark is anchored
        $bmNow = $d.Bookmarks.Item("_GoBack")
        $zRange = $d.Range($bmNow.Start, $bmNow.Start + 1)
        $zRange.Text = ""

        break
    }
}
